# Generate Report for Handoff
#
# The underlying localization-status report was regenerated: the
# "81ddff20-c635-4825-9f13-3fbb5b8226c2.md" file moved from "Handed back:
# in sync with en-US" into a fresh "Ready for handoff" state (new handoff
# file/datetime), and the three tracked files were re-sorted in the report
# (81ddff20... now sorts last, after ffffb80e2993... and ffffff6bbc32c2...).
#
# This script rewrites every data cell on all three sheets (Overview,
# zh-cn, de-de) to their final, post-regeneration values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value2 = "ffffb80e2993-1aae-4258-b44b-7f85ee356543.md"
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"

$ws.Range("A3").Value2 = "ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "Handed back: in sync with en-US"

$ws.Range("A4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "Ready for handoff"

$ws.Range("A5").Value2 = ".localization-config"
$ws.Range("B5").Value2 = "Not to be localized"
$ws.Range("C5").Value2 = "Not to be localized"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value2 = "ffffb80e2993-1aae-4258-b44b-7f85ee356543.md"
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf"
$ws.Range("D2").Value2 = "2016-02-22 17:56:45"
$ws.Range("E2").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.md"
$ws.Range("F2").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf"
$ws.Range("G2").Value2 = "2016-02-22 17:57:26"
$ws.Range("H2").Value2 = "Include"

$ws.Range("A3").Value2 = "ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf"
$ws.Range("D3").Value2 = "2016-02-22 17:56:45"
$ws.Range("E3").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.md"
$ws.Range("F3").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.zh-cn.xlf"
$ws.Range("G3").Value2 = "2016-02-22 17:57:26"
$ws.Range("H3").Value2 = "Include"

$ws.Range("A4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.zh-cn.xlf"
$ws.Range("D4").Value2 = "2016-02-22 18:00:43"
$ws.Range("E4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.md"
$ws.Range("F4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.zh-cn.xlf"
$ws.Range("G4").Value2 = "2016-02-22 17:59:35"
$ws.Range("H4").Value2 = "Include"

$ws.Range("A5").Value2 = ".localization-config"
$ws.Range("B5").Value2 = "Not to be localized"
$ws.Range("D5").Value2 = "0001-01-01 00:00:00"
$ws.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws.Range("H5").Value2 = "Ignored"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value2 = "ffffb80e2993-1aae-4258-b44b-7f85ee356543.md"
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf"
$ws.Range("D2").Value2 = "2016-02-22 17:56:56"
$ws.Range("E2").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.md"
$ws.Range("F2").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf"
$ws.Range("G2").Value2 = "2016-02-22 17:57:45"
$ws.Range("H2").Value2 = "Include"

$ws.Range("A3").Value2 = "ffffff6bbc32c2-1abb-4b71-9329-6a0c0e3a9fcd.md"
$ws.Range("B3").Value2 = "Handed back: in sync with en-US"
$ws.Range("C3").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf"
$ws.Range("D3").Value2 = "2016-02-22 17:56:56"
$ws.Range("E3").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.md"
$ws.Range("F3").Value2 = "f13d269d-2389-4baf-9322-ab170051d945.b476a835270279f96d7c74b645f3371bdbfad9e9.de-de.xlf"
$ws.Range("G3").Value2 = "2016-02-22 17:57:45"
$ws.Range("H3").Value2 = "Include"

$ws.Range("A4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.de-de.xlf"
$ws.Range("D4").Value2 = "2016-02-22 18:01:02"
$ws.Range("E4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.md"
$ws.Range("F4").Value2 = "81ddff20-c635-4825-9f13-3fbb5b8226c2.3f9735d69be03f5f557d0653e7a0dba0ef43cfa4.de-de.xlf"
$ws.Range("G4").Value2 = "2016-02-22 17:59:54"
$ws.Range("H4").Value2 = "Include"

$ws.Range("A5").Value2 = ".localization-config"
$ws.Range("B5").Value2 = "Not to be localized"
$ws.Range("D5").Value2 = "0001-01-01 00:00:00"
$ws.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws.Range("H5").Value2 = "Ignored"
